$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 previously had no SearchText (column A) value; add it now.
$ws.Range("A4").Value = "biology"

# New row 5 holding just a SearchText value.
$ws.Range("A5").Value = "tissue"

# Move the active selection to the newly added cell, matching the
# author's final cursor position in the saved workbook.
[void]$ws.Range("A5").Select()
